# Updated symbol list with refreshed price/volume/hour data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.21"
$ws.Range("E2").Value = "'-0.61%"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'27.35"
$ws.Range("E3").Value = "'3.92%"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.153"
$ws.Range("E4").Value = "'1.13%"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'0.05641"
$ws.Range("E5").Value = "'0.84%"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'6.467"
$ws.Range("E6").Value = "'-0.16%"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'0.8158"
$ws.Range("E7").Value = "'0.42%"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'0.8315"
$ws.Range("E8").Value = "'-1.16%"
$ws.Range("G8").Value = "'12"
$ws.Range("E9").Value = "'-1.41%"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.06923"
$ws.Range("E10").Value = "'-1.29%"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.02929"
$ws.Range("E11").Value = "'3.28%"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.09394"
$ws.Range("E12").Value = "'-0.13%"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.001509"
$ws.Range("E13").Value = "'-1.07%"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.04254"
$ws.Range("E14").Value = "'-9.74%"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'0.0005967"
$ws.Range("E15").Value = "'-93.92%"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.006133"
$ws.Range("E16").Value = "'0.41%"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'3.507"
$ws.Range("E17").Value = "'-1.58%"
$ws.Range("G17").Value = "'12"
$ws.Range("D18").Value = "'3.004"
$ws.Range("E18").Value = "'-0.96%"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'2.310"
$ws.Range("E19").Value = "'9.04%"
$ws.Range("G19").Value = "'12"
$ws.Range("E20").Value = "'-2.19%"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.03108"
$ws.Range("E21").Value = "'-3.59%"
$ws.Range("G21").Value = "'12"
$ws.Range("E22").Value = "'-2.17%"
$ws.Range("G22").Value = "'12"
$ws.Range("D23").Value = "'3.737"
$ws.Range("E23").Value = "'-0.10%"
$ws.Range("G23").Value = "'12"
$ws.Range("E24").Value = "'-0.12%"
$ws.Range("G24").Value = "'12"
$ws.Range("E25").Value = "'-1.90%"
$ws.Range("G25").Value = "'12"
$ws.Range("E26").Value = "'-2.80%"
$ws.Range("G26").Value = "'12"
$ws.Range("D27").Value = "'0.00009795"
$ws.Range("E27").Value = "'2.04%"
$ws.Range("G27").Value = "'12"
$ws.Range("E28").Value = "'-0.50%"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03644"
$ws.Range("E40").Value = "'-0.21%"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.006057"
$ws.Range("E41").Value = "'-0.53%"
$ws.Range("G41").Value = "'12"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.001799"
$ws.Range("E43").Value = "'-28.03%"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.008215"
$ws.Range("E44").Value = "'-4.09%"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005362"
$ws.Range("E45").Value = "'1.36%"
$ws.Range("G45").Value = "'12"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("G46").Value = "'12"
$ws.Range("D47").Value = "'0.1009"
$ws.Range("E47").Value = "'-23.98%"
$ws.Range("G47").Value = "'12"
$ws.Range("D48").Value = "'0.002652"
$ws.Range("E48").Value = "'29.22%"
$ws.Range("G48").Value = "'12"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("G49").Value = "'12"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
